# Applies the "symbol list" refresh for Mon Dec 12 17:06:49 UTC 2022.
# Every data row's "Hora" (G) cell advances 16 -> 17, several "Price" (D)
# cells get refreshed quotes, rows 6/7 and 42/43 swap places (their B/C/D/E
# cells trade values), and a couple of "Volume(1h)" (E) labels pick up a
# "Bestin24h"/"Worstin24h" suffix change. All of these values are plain text
# in the sheet (not real numbers), so each numeric-looking one is written with
# a temporary "@" (text) number format and then reset to the "Normal" style —
# this keeps Excel from "helpfully" converting e.g. "0.04690" into the number
# 0.0469 (dropping the trailing zero) or turning integers into floats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "276.09"  # D2
Set-TextValue 2 7 "17"  # G2

Set-TextValue 3 4 "21.01"  # D3
Set-TextValue 3 7 "17"  # G3

Set-TextValue 4 4 "6.231"  # D4
Set-TextValue 4 7 "17"  # G4

Set-TextValue 5 4 "0.06188"  # D5
Set-TextValue 5 7 "17"  # G5

$ws.Cells.Item(6, 2).Value = "GateToken"  # B6
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"  # C6
Set-TextValue 6 4 "3.579"  # D6
$ws.Cells.Item(6, 5).Value = "5GateTokenGT"  # E6
Set-TextValue 6 7 "17"  # G6

$ws.Cells.Item(7, 2).Value = "FTXToken"  # B7
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"  # C7
Set-TextValue 7 4 "1.543"  # D7
$ws.Cells.Item(7, 5).Value = "6FTXTokenFTT"  # E7
Set-TextValue 7 7 "17"  # G7

Set-TextValue 8 4 "6.567"  # D8
Set-TextValue 8 7 "17"  # G8

Set-TextValue 9 4 "0.8223"  # D9
Set-TextValue 9 7 "17"  # G9

Set-TextValue 10 7 "17"  # G10

Set-TextValue 11 4 "0.08211"  # D11
Set-TextValue 11 7 "17"  # G11

Set-TextValue 12 7 "17"  # G12

Set-TextValue 13 4 "0.03108"  # D13
Set-TextValue 13 7 "17"  # G13

Set-TextValue 14 4 "0.09126"  # D14
Set-TextValue 14 7 "17"  # G14

Set-TextValue 15 4 "3.773"  # D15
Set-TextValue 15 7 "17"  # G15

Set-TextValue 16 4 "0.001613"  # D16
Set-TextValue 16 7 "17"  # G16

Set-TextValue 17 4 "0.04690"  # D17
Set-TextValue 17 7 "17"  # G17

Set-TextValue 18 4 "0.006277"  # D18
Set-TextValue 18 7 "17"  # G18

Set-TextValue 19 4 "0.006141"  # D19
$ws.Cells.Item(19, 5).Value = "18HotbitTokenHTBBestin24h"  # E19
Set-TextValue 19 7 "17"  # G19

Set-TextValue 20 4 "0.001068"  # D20
Set-TextValue 20 7 "17"  # G20

Set-TextValue 21 7 "17"  # G21

Set-TextValue 22 4 "3.759"  # D22
Set-TextValue 22 7 "17"  # G22

Set-TextValue 23 4 "2.321"  # D23
Set-TextValue 23 7 "17"  # G23

Set-TextValue 24 4 "0.01385"  # D24
Set-TextValue 24 7 "17"  # G24

Set-TextValue 25 4 "0.3287"  # D25
Set-TextValue 25 7 "17"  # G25

Set-TextValue 26 7 "17"  # G26

Set-TextValue 27 7 "17"  # G27

Set-TextValue 28 4 "0.0002738"  # D28
Set-TextValue 28 7 "17"  # G28

Set-TextValue 29 7 "17"  # G29

Set-TextValue 30 7 "17"  # G30

Set-TextValue 31 7 "17"  # G31

Set-TextValue 32 7 "17"  # G32

Set-TextValue 33 7 "17"  # G33

Set-TextValue 34 7 "17"  # G34

Set-TextValue 35 7 "17"  # G35

Set-TextValue 36 7 "17"  # G36

Set-TextValue 37 7 "17"  # G37

Set-TextValue 38 7 "17"  # G38

Set-TextValue 39 7 "17"  # G39

Set-TextValue 40 4 "0.04682"  # D40
Set-TextValue 40 7 "17"  # G40

Set-TextValue 41 4 "0.007017"  # D41
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"  # E41
Set-TextValue 41 7 "17"  # G41

$ws.Cells.Item(42, 2).Value = "BKEXToken"  # B42
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"  # C42
Set-TextValue 42 4 "0.1105"  # D42
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"  # E42
Set-TextValue 42 7 "17"  # G42

$ws.Cells.Item(43, 2).Value = "CEJI"  # B43
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"  # C43
Set-TextValue 43 4 "0.003521"  # D43
$ws.Cells.Item(43, 5).Value = "42CEJICEJI"  # E43
Set-TextValue 43 7 "17"  # G43

Set-TextValue 44 4 "0.01114"  # D44
Set-TextValue 44 7 "17"  # G44

Set-TextValue 45 4 "0.00006263"  # D45
Set-TextValue 45 7 "17"  # G45

Set-TextValue 46 7 "17"  # G46

Set-TextValue 47 4 "0.8456"  # D47
Set-TextValue 47 7 "17"  # G47

Set-TextValue 48 7 "17"  # G48

Set-TextValue 49 4 "0.00001901"  # D49
Set-TextValue 49 7 "17"  # G49

Set-TextValue 50 4 "0.01241"  # D50
Set-TextValue 50 7 "17"  # G50

Set-TextValue 51 7 "17"  # G51
